$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GL_Date (C2): literal text "08/04/25" - force text so Excel doesn't
# auto-convert the date-looking string into a real date serial, then
# restore the default "Normal" style so no stray per-cell format sticks.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "08/04/25"
$ws.Range("C2").Style = "Normal"

# Vendor_Code (E2): plain text, no coercion risk.
$ws.Range("E2").Value = "LORSON"

# Invoice_Number (G2): was a plain number (106271), now becomes the text
# "01-871062" (contains a dash, so it must stay literal text).
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "01-871062"
$ws.Range("G2").Style = "Normal"

# Invoice_Date (I2): same literal-date-text situation as C2.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "08/04/25"
$ws.Range("I2").Style = "Normal"

# Invoice_Amount (J2): numeric value, straightforward.
$ws.Range("J2").Value = 89

# Job_Number (U2): numeric-looking text "25.09" must stay literal text.
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "25.09"
$ws.Range("U2").Style = "Normal"
